$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object "object[,]" 96,4
$arr[0,0] = 44393
$arr[0,1] = 12524
$arr[0,2] = 2454
$arr[0,3] = 14978
$arr[1,0] = 44394
$arr[1,1] = 12524
$arr[1,2] = 2454
$arr[1,3] = 14978
$arr[2,0] = 44395
$arr[2,1] = 12527
$arr[2,2] = 2456
$arr[2,3] = 14983
$arr[3,0] = 44396
$arr[3,1] = 12531
$arr[3,2] = 2457
$arr[3,3] = 14988
$arr[4,0] = 44397
$arr[4,1] = 12534
$arr[4,2] = 2457
$arr[4,3] = 14991
$arr[5,0] = 44398
$arr[5,1] = 12534
$arr[5,2] = 2457
$arr[5,3] = 14991
$arr[6,0] = 44399
$arr[6,1] = 12534
$arr[6,2] = 2457
$arr[6,3] = 14991
$arr[7,0] = 44400
$arr[7,1] = 12534
$arr[7,2] = 2457
$arr[7,3] = 14991
$arr[8,0] = 44401
$arr[8,1] = 12534
$arr[8,2] = 2457
$arr[8,3] = 14991
$arr[9,0] = 44402
$arr[9,1] = 12534
$arr[9,2] = 2457
$arr[9,3] = 14991
$arr[10,0] = 44403
$arr[10,1] = 12534
$arr[10,2] = 2457
$arr[10,3] = 14991
$arr[11,0] = 44404
$arr[11,1] = 12534
$arr[11,2] = 2457
$arr[11,3] = 14991
$arr[12,0] = 44405
$arr[12,1] = 12536
$arr[12,2] = 2457
$arr[12,3] = 14993
$arr[13,0] = 44406
$arr[13,1] = 12536
$arr[13,2] = 2457
$arr[13,3] = 14993
$arr[14,0] = 44407
$arr[14,1] = 12540
$arr[14,2] = 2457
$arr[14,3] = 14997
$arr[15,0] = 44408
$arr[15,1] = 12540
$arr[15,2] = 2457
$arr[15,3] = 14997
$arr[16,0] = 44409
$arr[16,1] = 12540
$arr[16,2] = 2457
$arr[16,3] = 14997
$arr[17,0] = 44410
$arr[17,1] = 12541
$arr[17,2] = 2457
$arr[17,3] = 14998
$arr[18,0] = 44411
$arr[18,1] = 12541
$arr[18,2] = 2457
$arr[18,3] = 14998
$arr[19,0] = 44412
$arr[19,1] = 12541
$arr[19,2] = 2457
$arr[19,3] = 14998
$arr[20,0] = 44413
$arr[20,1] = 12541
$arr[20,2] = 2457
$arr[20,3] = 14998
$arr[21,0] = 44414
$arr[21,1] = 12541
$arr[21,2] = 2457
$arr[21,3] = 14998
$arr[22,0] = 44415
$arr[22,1] = 12541
$arr[22,2] = 2457
$arr[22,3] = 14998
$arr[23,0] = 44416
$arr[23,1] = 12541
$arr[23,2] = 2457
$arr[23,3] = 14998
$arr[24,0] = 44417
$arr[24,1] = 12543
$arr[24,2] = 2457
$arr[24,3] = 15000
$arr[25,0] = 44418
$arr[25,1] = 12544
$arr[25,2] = 2457
$arr[25,3] = 15001
$arr[26,0] = 44419
$arr[26,1] = 12544
$arr[26,2] = 2457
$arr[26,3] = 15001
$arr[27,0] = 44420
$arr[27,1] = 12544
$arr[27,2] = 2457
$arr[27,3] = 15001
$arr[28,0] = 44421
$arr[28,1] = 12544
$arr[28,2] = 2457
$arr[28,3] = 15001
$arr[29,0] = 44422
$arr[29,1] = 12544
$arr[29,2] = 2457
$arr[29,3] = 15001
$arr[30,0] = 44423
$arr[30,1] = 12544
$arr[30,2] = 2457
$arr[30,3] = 15001
$arr[31,0] = 44424
$arr[31,1] = 12546
$arr[31,2] = 2457
$arr[31,3] = 15003
$arr[32,0] = 44425
$arr[32,1] = 12547
$arr[32,2] = 2457
$arr[32,3] = 15004
$arr[33,0] = 44426
$arr[33,1] = 12547
$arr[33,2] = 2457
$arr[33,3] = 15004
$arr[34,0] = 44427
$arr[34,1] = 12547
$arr[34,2] = 2457
$arr[34,3] = 15004
$arr[35,0] = 44428
$arr[35,1] = 12547
$arr[35,2] = 2457
$arr[35,3] = 15004
$arr[36,0] = 44429
$arr[36,1] = 12547
$arr[36,2] = 2457
$arr[36,3] = 15004
$arr[37,0] = 44430
$arr[37,1] = 12547
$arr[37,2] = 2457
$arr[37,3] = 15004
$arr[38,0] = 44431
$arr[38,1] = 12547
$arr[38,2] = 2457
$arr[38,3] = 15004
$arr[39,0] = 44432
$arr[39,1] = 12547
$arr[39,2] = 2457
$arr[39,3] = 15004
$arr[40,0] = 44433
$arr[40,1] = 12547
$arr[40,2] = 2457
$arr[40,3] = 15004
$arr[41,0] = 44434
$arr[41,1] = 12548
$arr[41,2] = 2457
$arr[41,3] = 15005
$arr[42,0] = 44435
$arr[42,1] = 12548
$arr[42,2] = 2457
$arr[42,3] = 15005
$arr[43,0] = 44436
$arr[43,1] = 12548
$arr[43,2] = 2457
$arr[43,3] = 15005
$arr[44,0] = 44437
$arr[44,1] = 12548
$arr[44,2] = 2457
$arr[44,3] = 15005
$arr[45,0] = 44438
$arr[45,1] = 12548
$arr[45,2] = 2457
$arr[45,3] = 15005
$arr[46,0] = 44439
$arr[46,1] = 12548
$arr[46,2] = 2457
$arr[46,3] = 15005
$arr[47,0] = 44440
$arr[47,1] = 12548
$arr[47,2] = 2457
$arr[47,3] = 15005
$arr[48,0] = 44441
$arr[48,1] = 12549
$arr[48,2] = 2457
$arr[48,3] = 15006
$arr[49,0] = 44442
$arr[49,1] = 12549
$arr[49,2] = 2457
$arr[49,3] = 15006
$arr[50,0] = 44443
$arr[50,1] = 12549
$arr[50,2] = 2457
$arr[50,3] = 15006
$arr[51,0] = 44444
$arr[51,1] = 12551
$arr[51,2] = 2458
$arr[51,3] = 15009
$arr[52,0] = 44445
$arr[52,1] = 12552
$arr[52,2] = 2458
$arr[52,3] = 15010
$arr[53,0] = 44446
$arr[53,1] = 12553
$arr[53,2] = 2458
$arr[53,3] = 15011
$arr[54,0] = 44447
$arr[54,1] = 12556
$arr[54,2] = 2458
$arr[54,3] = 15014
$arr[55,0] = 44448
$arr[55,1] = 12558
$arr[55,2] = 2458
$arr[55,3] = 15016
$arr[56,0] = 44449
$arr[56,1] = 12558
$arr[56,2] = 2458
$arr[56,3] = 15016
$arr[57,0] = 44450
$arr[57,1] = 12558
$arr[57,2] = 2458
$arr[57,3] = 15016
$arr[58,0] = 44451
$arr[58,1] = 12560
$arr[58,2] = 2458
$arr[58,3] = 15018
$arr[59,0] = 44452
$arr[59,1] = 12562
$arr[59,2] = 2458
$arr[59,3] = 15020
$arr[60,0] = 44453
$arr[60,1] = 12562
$arr[60,2] = 2458
$arr[60,3] = 15020
$arr[61,0] = 44454
$arr[61,1] = 12562
$arr[61,2] = 2458
$arr[61,3] = 15020
$arr[62,0] = 44455
$arr[62,1] = 12569
$arr[62,2] = 2459
$arr[62,3] = 15028
$arr[63,0] = 44456
$arr[63,1] = 12569
$arr[63,2] = 2459
$arr[63,3] = 15028
$arr[64,0] = 44457
$arr[64,1] = 12573
$arr[64,2] = 2461
$arr[64,3] = 15034
$arr[65,0] = 44458
$arr[65,1] = 12580
$arr[65,2] = 2461
$arr[65,3] = 15041
$arr[66,0] = 44459
$arr[66,1] = 12589
$arr[66,2] = 2461
$arr[66,3] = 15050
$arr[67,0] = 44460
$arr[67,1] = 12589
$arr[67,2] = 2461
$arr[67,3] = 15050
$arr[68,0] = 44461
$arr[68,1] = 12594
$arr[68,2] = 2461
$arr[68,3] = 15055
$arr[69,0] = 44462
$arr[69,1] = 12594
$arr[69,2] = 2461
$arr[69,3] = 15055
$arr[70,0] = 44463
$arr[70,1] = 12596
$arr[70,2] = 2462
$arr[70,3] = 15058
$arr[71,0] = 44464
$arr[71,1] = 12596
$arr[71,2] = 2462
$arr[71,3] = 15058
$arr[72,0] = 44465
$arr[72,1] = 12596
$arr[72,2] = 2462
$arr[72,3] = 15058
$arr[73,0] = 44466
$arr[73,1] = 12606
$arr[73,2] = 2462
$arr[73,3] = 15068
$arr[74,0] = 44467
$arr[74,1] = 12620
$arr[74,2] = 2462
$arr[74,3] = 15082
$arr[75,0] = 44468
$arr[75,1] = 12637
$arr[75,2] = 2468
$arr[75,3] = 15105
$arr[76,0] = 44469
$arr[76,1] = 12649
$arr[76,2] = 2472
$arr[76,3] = 15121
$arr[77,0] = 44470
$arr[77,1] = 12660
$arr[77,2] = 2476
$arr[77,3] = 15136
$arr[78,0] = 44471
$arr[78,1] = 12668
$arr[78,2] = 2476
$arr[78,3] = 15144
$arr[79,0] = 44472
$arr[79,1] = 12676
$arr[79,2] = 2476
$arr[79,3] = 15152
$arr[80,0] = 44473
$arr[80,1] = 12690
$arr[80,2] = 2479
$arr[80,3] = 15169
$arr[81,0] = 44474
$arr[81,1] = 12697
$arr[81,2] = 2479
$arr[81,3] = 15176
$arr[82,0] = 44475
$arr[82,1] = 12705
$arr[82,2] = 2479
$arr[82,3] = 15184
$arr[83,0] = 44476
$arr[83,1] = 12720
$arr[83,2] = 2482
$arr[83,3] = 15202
$arr[84,0] = 44477
$arr[84,1] = 12729
$arr[84,2] = 2488
$arr[84,3] = 15217
$arr[85,0] = 44478
$arr[85,1] = 12731
$arr[85,2] = 2488
$arr[85,3] = 15219
$arr[86,0] = 44479
$arr[86,1] = 12735
$arr[86,2] = 2489
$arr[86,3] = 15224
$arr[87,0] = 44480
$arr[87,1] = 12752
$arr[87,2] = 2499
$arr[87,3] = 15251
$arr[88,0] = 44481
$arr[88,1] = 12772
$arr[88,2] = 2501
$arr[88,3] = 15273
$arr[89,0] = 44482
$arr[89,1] = 12791
$arr[89,2] = 2508
$arr[89,3] = 15299
$arr[90,0] = 44483
$arr[90,1] = 12816
$arr[90,2] = 2508
$arr[90,3] = 15324
$arr[91,0] = 44484
$arr[91,1] = 12833
$arr[91,2] = 2508
$arr[91,3] = 15341
$arr[92,0] = 44485
$arr[92,1] = 12846
$arr[92,2] = 2508
$arr[92,3] = 15354
$arr[93,0] = 44486
$arr[93,1] = 12854
$arr[93,2] = 2508
$arr[93,3] = 15362
$arr[94,0] = 44487
$arr[94,1] = 12864
$arr[94,2] = 2509
$arr[94,3] = 15373
$arr[95,0] = 44488
$arr[95,1] = 12872
$arr[95,2] = 2509
$arr[95,3] = 15381
$rng = $ws.Range("A274:D369")
$rng.Value = $arr
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
